$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "63.105.66"
$ws.Range("E2").Value = "  +1.55%  "
Set-TextValue $ws.Range("D3") "2.470.35"
$ws.Range("E3").Value = "  +2.53%  "
$ws.Range("E4").Value = "  -0.25%  "
Set-TextValue $ws.Range("D5") "576.51"
$ws.Range("E5").Value = "  +0.85%  "
Set-TextValue $ws.Range("D6") "146.19"
$ws.Range("E6").Value = "  +1.35%  "
$ws.Range("E7").Value = "  +0.28%  "
Set-TextValue $ws.Range("D8") "0.541"
$ws.Range("E8").Value = "  +0.73%  "
Set-TextValue $ws.Range("D9") "2.469.80"
$ws.Range("E9").Value = "  +1.58%  "
Set-TextValue $ws.Range("D10") "0.111"
$ws.Range("E10").Value = "  +1.48%  "
$ws.Range("E11").Value = "  +1.30%  "
Set-TextValue $ws.Range("D12") "5.29"
$ws.Range("E12").Value = "  +0.77%  "
Set-TextValue $ws.Range("D13") "0.355"
$ws.Range("E13").Value = "  +1.74%  "
Set-TextValue $ws.Range("D14") "29.14"
$ws.Range("E14").Value = "  +9.27%  "
$ws.Range("E15").Value = "  +0.63%  "
Set-TextValue $ws.Range("D16") "2.918.61"
$ws.Range("E16").Value = "  +2.43%  "
Set-TextValue $ws.Range("D17") "63.158.78"
$ws.Range("E17").Value = "  +2.06%  "
Set-TextValue $ws.Range("D18") "2.466.27"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("E19").Value = "  +0.05%  "
Set-TextValue $ws.Range("D20") "11.10"
$ws.Range("E20").Value = "  +2.51%  "
Set-TextValue $ws.Range("D21") "330.33"
$ws.Range("E21").Value = "  +1.41%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("E23").Value = "  +9.29%  "
$ws.Range("E24").Value = "  +0.04%  "
Set-TextValue $ws.Range("D25") "66.25"
$ws.Range("E25").Value = "  +1.56%  "
Set-TextValue $ws.Range("D26") "663.58"
$ws.Range("E26").Value = "  +8.23%  "
Set-TextValue $ws.Range("D27") "9.21"
$ws.Range("E27").Value = "  +9.39%  "
Set-TextValue $ws.Range("D28") "0.0000100"
$ws.Range("E28").Value = "  +1.77%  "
Set-TextValue $ws.Range("D29") "2.603.27"
$ws.Range("E29").Value = "  +2.58%  "
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("E31").Value = "  +3.89%  "
Set-TextValue $ws.Range("D32") "8.15"
$ws.Range("E32").Value = "  +1.04%  "
Set-TextValue $ws.Range("D33") "1.87"
$ws.Range("E33").Value = "  +2.65%  "
Set-TextValue $ws.Range("D34") "0.137"
$ws.Range("E34").Value = "  +1.60%  "
$ws.Range("E35").Value = "  +4.39%  "
Set-TextValue $ws.Range("D36") "0.998"
$ws.Range("E37").Value = "  +0.79%  "
Set-TextValue $ws.Range("D38") "5.51"
$ws.Range("E38").Value = "  +2.31%  "
$ws.Range("E39").Value = "  +0.46%  "
Set-TextValue $ws.Range("D40") "153.00"
$ws.Range("E40").Value = "  +0.04%  "
Set-TextValue $ws.Range("D41") "18.83"
$ws.Range("E41").Value = "  +1.66%  "
Set-TextValue $ws.Range("D42") "2.73"
$ws.Range("E42").Value = "  +1.96%  "
$ws.Range("E43").Value = "  +1.85%  "
$ws.Range("E44").Value = "  +0.05%  "
Set-TextValue $ws.Range("D45") "0.0₆0299"
$ws.Range("E45").Value = "  +5.72%  "
$ws.Range("E46").Value = "  +27.48%  "
Set-TextValue $ws.Range("D47") "147.42"
$ws.Range("E47").Value = "  +3.12%  "
Set-TextValue $ws.Range("D48") "3.63"
$ws.Range("E48").Value = "  +1.47%  "
Set-TextValue $ws.Range("D49") "20.88"
$ws.Range("E49").Value = "  +2.90%  "
Set-TextValue $ws.Range("D50") "0.608"
$ws.Range("E50").Value = "  +1.79%  "
Set-TextValue $ws.Range("D51") "0.0517"
$ws.Range("E51").Value = "  +0.89%  "
